$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.993.05"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "3.718.37"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.80"
$ws.Range("E5").Value = "  +7.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.62"
$ws.Range("E6").Value = "  +10.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.717"
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.72"
$ws.Range("E11").Value = "  +8.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.56"
$ws.Range("E13").Value = "  -0.70%  "
$ws.Range("D14").Value = "4.305.22"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("D15").Value = "3.716.02"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.30"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.127"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.14"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "68.783.98"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "411.70"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.60"
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.52"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.05"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.85"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.93"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.04"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.77"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.68"
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.13"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.33"
$ws.Range("E31").Value = "  -8.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.71"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "45.57"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "631.40"
$ws.Range("E35").Value = "  +5.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.67"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.416"
$ws.Range("E37").Value = "  +4.14%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0820"
$ws.Range("E38").Value = "  -11.00%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  +3.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0444"
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.140"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("D46").Value = "2.874.49"
$ws.Range("E46").Value = "  +6.14%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.06"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.96"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.09"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.57"
$ws.Range("E51").Value = "  -18.54%  "
